$d = $word.ActiveDocument

# 1) "The Pokemon table ... Their primary key will be the pokemon name since they are
#    unique. It will also connect ..." -> "... primary key will be the id. It will also
#    connect ..."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The Pokemon table with have most of the start columns*") {
        $p.Range.Find.Execute("pokemon name since they are unique", $true, $false, $false, $false, $false, $true, 1, $false, "id", 2)
        break
    }
}

# 2) "The joiner table will have the name from the pokemon table ..." -> "... will have
#    the id from the pokemon table ..."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The joiner table will have the name*") {
        $p.Range.Find.Execute("have the name from", $true, $false, $false, $false, $false, $true, 1, $false, "have the id from", 2)
        break
    }
}

# 3) "Finally the type table will have the types, it only has one column as the types
#    are all unique so it is also the primary key. " -> "Finally the type table will
#    have the id as primary key and then the type."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Finally the type table will have the*") {
        $p.Range.Find.Execute("types, it only has one column as the types are all unique so it is also the primary key. ", $true, $false, $false, $false, $false, $true, 1, $false, "id as primary key and then the type.", 2)
        break
    }
}

# 4) "- name: string (varchar sql PK)" -> "- name: string (varchar sql)"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- name: string*") {
        $p.Range.Find.Execute(" PK)", $true, $false, $false, $false, $false, $true, 1, $false, ")", 2)
        break
    }
}

# 5) Insert a new paragraph "- id: int (PK)" right after "- number: int (FK)" (end of
#    the Pokemon model block, before the blank line preceding "Model: Pokedex ...").
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- number: int (FK)*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "- id: int (PK)"
        break
    }
}

# 6) "- pokemon_name: string (varchar sql FK)" -> "- pokemon_id: string (varchar sql FK)"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- pokemon_name:*") {
        $p.Range.Find.Execute("pokemon_name", $true, $false, $false, $false, $false, $true, 1, $false, "pokemon_id", 2)
        break
    }
}

# 7) "- type_type: string (varchar sql FK)" -> "- type_id: string (varchar sql FK)"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- type_type:*") {
        $p.Range.Find.Execute("type_type", $true, $false, $false, $false, $false, $true, 1, $false, "type_id", 2)
        break
    }
}

# 8) Insert a new paragraph "-id: int (PK)" right before "- type: string (varchar sql
#    PK)" in the "Types Table" block, then strip " PK" from that existing paragraph so
#    it reads "- type: string (varchar sql)".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Types Table (many to many*") {
        $nextPara = $p.Next()
        $nextPara.Range.InsertParagraphBefore()
        $newPara = $p.Next()
        $newPara.Range.Text = "-id: int (PK)"
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- type: string*") {
        $p.Range.Find.Execute(" PK)", $true, $false, $false, $false, $false, $true, 1, $false, ")", 2)
        break
    }
}
